$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the recruitment band elevation criteria values
$ws.Range("C23").Value = 0
$ws.Range("C24").Value = 150

# Update the selected cell / view state to C25
$ws.Range("C25").Select()
